$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Unit  Type" (note: two spaces) in column V, matching the
# existing header formatting by copying format from the adjacent header cell
$ws.Range("V1").Value2 = "Unit  Type"
$ws.Range("U1").Copy() | Out-Null
$ws.Range("V1").PasteSpecial(-4122) | Out-Null

# New data cell V2 mirrors H2's value ("box")
$ws.Range("V2").Value2 = $ws.Range("H2").Value2

# Give the new column a sensible custom width (matches the ~10 char width
# used for this field in the final workbook)
$ws.Columns.Item(22).ColumnWidth = 9.14

# Update selection to mirror the recorded end-state of the edit session
$ws.Range("G19").Select() | Out-Null
